# Update Name of Algo
# Apply updated imputed values in column E for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "E9"  = 17.346
    "E18" = 16.593
    "E20" = 16.314
    "E27" = 16.531
    "E35" = 16.407
    "E69" = 17.4
    "E76" = 16.433
    "E78" = 16.507
    "E82" = 16.779
    "E83" = 16.935
    "E93" = 17.103
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
